$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "24/7/2023"
$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "Online meeting with Classplus"

$ws.Columns.Item(3).ColumnWidth = 27.6

$ws.Range("D9").Select()
